$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").EntireColumn.Delete()

$ws.Range("E6").Select()
